# Apply the "added acceptable answers to cancer screening" edit:
# Replace the short placeholder "Answer Options" values in column C of the
# Cancer sheet with the fuller, more descriptive acceptable-answer lists.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Cancer")

$yesNo     = "Yes;No;Refused;Don't Know;Missing"
$urinary   = "Never;Less than once a month;A few times a month;A few times a week; Ever day and/or night; Refused;Don't Know;Missing"
$hearing   = "Excellent;Good;A little trouble;Moderate hearing trouble;A lot of trouble;Deaf;Refused;Don't Know;Missing"
$dentist   = "Went in on own for check-up, examination, or cleaning;Was called in by dentist for check-up, examination, or cleaning;Something was wrong, bothering, or hurting; Went for treatment of a condition that dentist previously discovered;Other;Refused;Don't Know;Missing"
$tiredness = "Not at all;Several Days;More than half the days;Nearly every day;Refused;Don't Know;Missing"

$ws.Range("C2").Value  = $yesNo      # INQ030  - Do You Get Income From Social Security?
$ws.Range("C3").Value  = $urinary    # KIQ005  - How Often Do You Have Urinary Leakage per week?
$ws.Range("C4").Value  = $hearing    # AUQ054  - General Condition of Hearing?
$ws.Range("C5").Value  = $yesNo      # MCQ560  - Have You Ever Had GallBladder Surgery?
$ws.Range("C6").Value  = $yesNo      # MCQ371D - Are you watching your weight?
$ws.Range("C7").Value  = $yesNo      # HUQ071  - Were you a Patient in Hospital Overnight?
$ws.Range("C8").Value  = $dentist    # OHQ033  - Main Reason for Visiting Dentist?
$ws.Range("C9").Value  = $yesNo      # MCQ092  - Received Blood Transfusion?
$ws.Range("C10").Value = $yesNo      # SMQ020  - Have You Smoked Atleast 100 Cigarettes?
$ws.Range("C11").Value = $tiredness  # DPQ040  - Over the Last Two Weeks have You Felt Tired...
$ws.Range("C12").Value = $yesNo      # PUQ110  - Any Chemical Products Used to Kill Weeds?

# Widen column C so the longer answer-option text fits (matches the
# "best fit" auto-sizing Excel performs after the content changes).
$ws.Columns.Item(3).ColumnWidth = 233.5

# The user had moved their selection to C20 in the Cancer sheet before saving.
$ws.Range("C20").Select() | Out-Null
